$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (column C) and montant_total (column D) for rows with new data
# as of 2020-08-14 refresh
$ws.Cells.Item(2, 3).Value = 317195
$ws.Cells.Item(2, 4).Value = 404245959
$ws.Cells.Item(8, 3).Value = 854
$ws.Cells.Item(8, 4).Value = 1256408
$ws.Cells.Item(10, 3).Value = 116261
$ws.Cells.Item(10, 4).Value = 170357721
$ws.Cells.Item(12, 3).Value = 58767
$ws.Cells.Item(12, 4).Value = 84807331
$ws.Cells.Item(16, 3).Value = 3989
$ws.Cells.Item(16, 4).Value = 5660873
$ws.Cells.Item(20, 3).Value = 6529
$ws.Cells.Item(20, 4).Value = 9106434
$ws.Cells.Item(22, 3).Value = 76621
$ws.Cells.Item(22, 4).Value = 95578711
$ws.Cells.Item(28, 3).Value = 32256
$ws.Cells.Item(28, 4).Value = 47220910
$ws.Cells.Item(30, 3).Value = 11379
$ws.Cells.Item(30, 4).Value = 16366039
$ws.Cells.Item(33, 3).Value = 1558
$ws.Cells.Item(33, 4).Value = 2188307
$ws.Cells.Item(35, 3).Value = 1792
$ws.Cells.Item(35, 4).Value = 2528531
$ws.Cells.Item(36, 3).Value = 96264
$ws.Cells.Item(36, 4).Value = 121206635
$ws.Cells.Item(44, 3).Value = 44114
$ws.Cells.Item(44, 4).Value = 64648177
$ws.Cells.Item(46, 3).Value = 9051
$ws.Cells.Item(46, 4).Value = 12989169
$ws.Cells.Item(48, 3).Value = 1399
$ws.Cells.Item(48, 4).Value = 1943109
$ws.Cells.Item(51, 3).Value = 2265
$ws.Cells.Item(51, 4).Value = 3159852
$ws.Cells.Item(52, 3).Value = 68374
$ws.Cells.Item(52, 4).Value = 85788682
$ws.Cells.Item(58, 3).Value = 27951
$ws.Cells.Item(58, 4).Value = 40992109
$ws.Cells.Item(61, 3).Value = 10983
$ws.Cells.Item(61, 4).Value = 15881241
$ws.Cells.Item(67, 3).Value = 1448
$ws.Cells.Item(67, 4).Value = 2028065
$ws.Cells.Item(69, 3).Value = 20296
$ws.Cells.Item(69, 4).Value = 26583360
$ws.Cells.Item(73, 3).Value = 7533
$ws.Cells.Item(73, 4).Value = 11029021
$ws.Cells.Item(75, 3).Value = 5071
$ws.Cells.Item(75, 4).Value = 7362678
$ws.Cells.Item(78, 3).Value = 139331
$ws.Cells.Item(78, 4).Value = 173759993
$ws.Cells.Item(84, 3).Value = 63094
$ws.Cells.Item(84, 4).Value = 92474537
$ws.Cells.Item(87, 3).Value = 29429
$ws.Cells.Item(87, 4).Value = 42568516
$ws.Cells.Item(90, 3).Value = 2774
$ws.Cells.Item(90, 4).Value = 3917945
$ws.Cells.Item(91, 3).Value = 32464
$ws.Cells.Item(91, 4).Value = 43975596
$ws.Cells.Item(95, 3).Value = 7852
$ws.Cells.Item(95, 4).Value = 11545581
$ws.Cells.Item(97, 3).Value = 7173
$ws.Cells.Item(97, 4).Value = 10400082
$ws.Cells.Item(99, 3).Value = 528
$ws.Cells.Item(99, 4).Value = 750705
$ws.Cells.Item(100, 3).Value = 484
$ws.Cells.Item(100, 4).Value = 698443
$ws.Cells.Item(101, 3).Value = 8947
$ws.Cells.Item(101, 4).Value = 12421784
$ws.Cells.Item(103, 3).Value = 2249
$ws.Cells.Item(103, 4).Value = 3313802
$ws.Cells.Item(105, 3).Value = 3015
$ws.Cells.Item(105, 4).Value = 4402834
$ws.Cells.Item(109, 3).Value = 139790
$ws.Cells.Item(109, 4).Value = 172845536
$ws.Cells.Item(115, 3).Value = 52319
$ws.Cells.Item(115, 4).Value = 76694938
$ws.Cells.Item(117, 3).Value = 26694
$ws.Cells.Item(117, 4).Value = 38672513
$ws.Cells.Item(121, 3).Value = 2213
$ws.Cells.Item(121, 4).Value = 3108355
$ws.Cells.Item(123, 3).Value = 496217
$ws.Cells.Item(123, 4).Value = 654354346
$ws.Cells.Item(130, 3).Value = 205487
$ws.Cells.Item(130, 4).Value = 302063868
$ws.Cells.Item(133, 3).Value = 177731
$ws.Cells.Item(133, 4).Value = 258330503
$ws.Cells.Item(136, 3).Value = 2835
$ws.Cells.Item(136, 4).Value = 3985381
$ws.Cells.Item(138, 3).Value = 6222
$ws.Cells.Item(138, 4).Value = 8791173
$ws.Cells.Item(141, 3).Value = 44006
$ws.Cells.Item(141, 4).Value = 58741729
$ws.Cells.Item(147, 3).Value = 13939
$ws.Cells.Item(147, 4).Value = 20440478
$ws.Cells.Item(148, 3).Value = 3712
$ws.Cells.Item(148, 4).Value = 5352680
$ws.Cells.Item(154, 3).Value = 17363
$ws.Cells.Item(154, 4).Value = 22940196
$ws.Cells.Item(157, 3).Value = 54
$ws.Cells.Item(157, 4).Value = 79406
$ws.Cells.Item(158, 3).Value = 7089
$ws.Cells.Item(158, 4).Value = 10309713
$ws.Cells.Item(160, 3).Value = 4936
$ws.Cells.Item(160, 4).Value = 7103063
$ws.Cells.Item(163, 3).Value = 265
$ws.Cells.Item(163, 4).Value = 378864
$ws.Cells.Item(165, 3).Value = 15685
$ws.Cells.Item(165, 4).Value = 22757845
$ws.Cells.Item(166, 3).Value = 1765
$ws.Cells.Item(166, 4).Value = 2625230
$ws.Cells.Item(170, 3).Value = 83
$ws.Cells.Item(170, 4).Value = 124449
$ws.Cells.Item(171, 3).Value = 86813
$ws.Cells.Item(171, 4).Value = 108585519
$ws.Cells.Item(178, 3).Value = 33611
$ws.Cells.Item(178, 4).Value = 49289821
$ws.Cells.Item(180, 3).Value = 12871
$ws.Cells.Item(180, 4).Value = 18595409
$ws.Cells.Item(182, 3).Value = 1242
$ws.Cells.Item(182, 4).Value = 1738896
$ws.Cells.Item(184, 3).Value = 1620
$ws.Cells.Item(184, 4).Value = 2275662
$ws.Cells.Item(186, 3).Value = 236145
$ws.Cells.Item(186, 4).Value = 293546188
$ws.Cells.Item(194, 3).Value = 86018
$ws.Cells.Item(194, 4).Value = 126089169
$ws.Cells.Item(197, 3).Value = 32720
$ws.Cells.Item(197, 4).Value = 47088886
$ws.Cells.Item(200, 3).Value = 5089
$ws.Cells.Item(200, 4).Value = 7248543
$ws.Cells.Item(203, 3).Value = 4793
$ws.Cells.Item(203, 4).Value = 6634058
$ws.Cells.Item(206, 3).Value = 261090
$ws.Cells.Item(206, 4).Value = 323122381
$ws.Cells.Item(208, 3).Value = 253
$ws.Cells.Item(208, 4).Value = 362087
$ws.Cells.Item(215, 3).Value = 94475
$ws.Cells.Item(215, 4).Value = 138207104
$ws.Cells.Item(218, 3).Value = 50929
$ws.Cells.Item(218, 4).Value = 73601651
$ws.Cells.Item(221, 3).Value = 4654
$ws.Cells.Item(221, 4).Value = 6534091
$ws.Cells.Item(224, 3).Value = 5644
$ws.Cells.Item(224, 4).Value = 7802387
$ws.Cells.Item(227, 3).Value = 105097
$ws.Cells.Item(227, 4).Value = 131481761
$ws.Cells.Item(232, 3).Value = 563
$ws.Cells.Item(232, 4).Value = 822439
$ws.Cells.Item(234, 3).Value = 49141
$ws.Cells.Item(234, 4).Value = 71991405
$ws.Cells.Item(236, 3).Value = 12249
$ws.Cells.Item(236, 4).Value = 17610577
$ws.Cells.Item(238, 3).Value = 1888
$ws.Cells.Item(238, 4).Value = 2705882
$ws.Cells.Item(240, 3).Value = 2463
$ws.Cells.Item(240, 4).Value = 3442565
$ws.Cells.Item(241, 3).Value = 254580
$ws.Cells.Item(241, 4).Value = 321414548
$ws.Cells.Item(247, 3).Value = 821
$ws.Cells.Item(247, 4).Value = 1205563
$ws.Cells.Item(249, 3).Value = 95033
$ws.Cells.Item(249, 4).Value = 139245527
$ws.Cells.Item(252, 3).Value = 64227
$ws.Cells.Item(252, 4).Value = 93066999
$ws.Cells.Item(254, 3).Value = 2398
$ws.Cells.Item(254, 4).Value = 3383828
$ws.Cells.Item(257, 3).Value = 4526
$ws.Cells.Item(257, 4).Value = 6354616
